$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# hunk 0: sheet ALC
$ws.Cells.Item(6, 8).Value = 10755891
$ws.Cells.Item(6, 9).Value = 22228336
$ws.Cells.Item(6, 11).Value = 66685008
$ws.Cells.Item(6, 13).Value = -66684896

# hunk 1: sheet ALC
$ws.Cells.Item(32, 8).Value = 5874.5
$ws.Cells.Item(32, 10).Value = 6499.6665
$ws.Cells.Item(32, 12).Value = 6499.6665
$ws.Cells.Item(32, 14).Value = -7151.6665

# hunk 2: sheet ALC
$ws.Cells.Item(98, 8).Value = 29556.87
$ws.Cells.Item(98, 9).Value = 35113.375
$ws.Cells.Item(98, 11).Value = 35113.375
$ws.Cells.Item(98, 13).Value = -33615.375

# hunk 3: sheet ALC
$ws.Cells.Item(122, 8).Value = 29556.87
$ws.Cells.Item(122, 9).Value = 35113.375
$ws.Cells.Item(122, 11).Value = 105340.125
$ws.Cells.Item(122, 13).Value = -102890.125

# hunk 4: sheet ALC
$ws.Cells.Item(131, 8).Value = 4096
$ws.Cells.Item(131, 10).Value = 7558.5
$ws.Cells.Item(131, 12).Value = 22675.5
$ws.Cells.Item(131, 14).Value = -32755.5

# hunk 5: sheet ALC
$ws.Cells.Item(132, 8).Value = 3590.2122
$ws.Cells.Item(132, 9).Value = 3467.0322
$ws.Cells.Item(132, 10).Value = 5499.5
$ws.Cells.Item(132, 11).Value = 10401.0966
$ws.Cells.Item(132, 12).Value = 16498.5
$ws.Cells.Item(132, 13).Value = -7871.096600000001
$ws.Cells.Item(132, 14).Value = -21558.5

# hunk 6: sheet ALC
$ws.Cells.Item(138, 8).Value = 3700.7795
$ws.Cells.Item(138, 9).Value = 606.8889
$ws.Cells.Item(138, 10).Value = 4257.68
$ws.Cells.Item(138, 11).Value = 1820.6667
$ws.Cells.Item(138, 12).Value = 12773.04
$ws.Cells.Item(138, 13).Value = 3319.3333
$ws.Cells.Item(138, 14).Value = -23053.04

$ws = $wb.Worksheets.Item("ARM")
# hunk 7: sheet ARM
$ws.Cells.Item(32, 8).Value = 1908.8868
$ws.Cells.Item(32, 9).Value = 1975.2549
$ws.Cells.Item(32, 11).Value = 1975.2549
$ws.Cells.Item(32, 13).Value = -1688.2549

# hunk 8: sheet ARM
$ws.Cells.Item(45, 8).Value = 5751.067
$ws.Cells.Item(45, 9).Value = 5183.5
$ws.Cells.Item(45, 11).Value = 5183.5
$ws.Cells.Item(45, 13).Value = -4806.5

# hunk 9: sheet ARM
$ws.Cells.Item(122, 8).Value = 1756337.2
$ws.Cells.Item(122, 9).Value = 5174.75
$ws.Cells.Item(122, 10).Value = 3507499.8
$ws.Cells.Item(122, 11).Value = 15524.25
$ws.Cells.Item(122, 12).Value = 10522499.4
$ws.Cells.Item(122, 13).Value = -13074.25
$ws.Cells.Item(122, 14).Value = -10527399.4

$ws = $wb.Worksheets.Item("BSM")
# hunk 10: sheet BSM
$ws.Cells.Item(134, 8).Value = 4442.524
$ws.Cells.Item(134, 9).Value = 3173.0667
$ws.Cells.Item(134, 11).Value = 9519.2001
$ws.Cells.Item(134, 13).Value = -6984.2001

$ws = $wb.Worksheets.Item("CRP")
# hunk 11: sheet CRP
$ws.Cells.Item(22, 8).Value = 666.4706
$ws.Cells.Item(22, 9).Value = 601
$ws.Cells.Item(22, 10).Value = 670.5625
$ws.Cells.Item(22, 11).Value = 601
$ws.Cells.Item(22, 12).Value = 670.5625
$ws.Cells.Item(22, 13).Value = -251
$ws.Cells.Item(22, 14).Value = -1370.5625

$ws = $wb.Worksheets.Item("CUL")
# hunk 12: sheet CUL
$ws.Cells.Item(80, 8).Value = 105099.75
$ws.Cells.Item(80, 9).Value = 1449.5
$ws.Cells.Item(80, 10).Value = 139649.83
$ws.Cells.Item(80, 11).Value = 4348.5
$ws.Cells.Item(80, 12).Value = 418949.49
$ws.Cells.Item(80, 13).Value = -3412.5
$ws.Cells.Item(80, 14).Value = -420821.49

# hunk 13: sheet CUL
$ws.Cells.Item(83, 8).Value = 105099.75
$ws.Cells.Item(83, 9).Value = 1449.5
$ws.Cells.Item(83, 10).Value = 139649.83
$ws.Cells.Item(83, 11).Value = 13045.5
$ws.Cells.Item(83, 12).Value = 1256848.47
$ws.Cells.Item(83, 13).Value = -8365.5
$ws.Cells.Item(83, 14).Value = -1266208.47

# hunk 14: sheet CUL
$ws.Cells.Item(131, 8).Value = 47620384
$ws.Cells.Item(131, 9).Value = 90909740
$ws.Cells.Item(131, 11).Value = 272729220
$ws.Cells.Item(131, 13).Value = -272724180

# hunk 15: sheet CUL
$ws.Cells.Item(137, 8).Value = 2788
$ws.Cells.Item(137, 9).Value = 3142.8572
$ws.Cells.Item(137, 11).Value = 9428.571599999999
$ws.Cells.Item(137, 13).Value = -4328.571599999999

$ws = $wb.Worksheets.Item("GSM")
# hunk 16: sheet GSM
$ws.Cells.Item(2, 8).Value = 1400.7693
$ws.Cells.Item(2, 9).Value = 1564.5454
$ws.Cells.Item(2, 11).Value = 1564.5454
$ws.Cells.Item(2, 13).Value = -1451.5454

# hunk 17: sheet GSM
$ws.Cells.Item(97, 8).Value = 7565.697
$ws.Cells.Item(97, 9).Value = 8376.074000000001
$ws.Cells.Item(97, 10).Value = 3919
$ws.Cells.Item(97, 11).Value = 8376.074000000001
$ws.Cells.Item(97, 12).Value = 3919
$ws.Cells.Item(97, 13).Value = -7880.074000000001
$ws.Cells.Item(97, 14).Value = -4911

# hunk 18: sheet GSM
$ws.Cells.Item(122, 8).Value = 35699.6
$ws.Cells.Item(122, 9).Value = 43999.5
$ws.Cells.Item(122, 10).Value = 30166.334
$ws.Cells.Item(122, 11).Value = 131998.5
$ws.Cells.Item(122, 12).Value = 90499.00199999999
$ws.Cells.Item(122, 13).Value = -129548.5
$ws.Cells.Item(122, 14).Value = -95399.00199999999

$ws = $wb.Worksheets.Item("LTW")
# hunk 19: sheet LTW
$ws.Cells.Item(7, 8).Value = 41372.332
$ws.Cells.Item(7, 9).Value = 63853.145
$ws.Cells.Item(7, 11).Value = 63853.145
$ws.Cells.Item(7, 13).Value = -63741.145

# hunk 20: sheet LTW
$ws.Cells.Item(40, 8).Value = 99796.60000000001
$ws.Cells.Item(40, 9).Value = 201999.5
$ws.Cells.Item(40, 10).Value = 31661.334
$ws.Cells.Item(40, 11).Value = 201999.5
$ws.Cells.Item(40, 12).Value = 31661.334
$ws.Cells.Item(40, 13).Value = -201863.5
$ws.Cells.Item(40, 14).Value = -31933.334

# hunk 21: sheet LTW
$ws.Cells.Item(122, 8).Value = 3262.8667
$ws.Cells.Item(122, 9).Value = 2995.6924
$ws.Cells.Item(122, 10).Value = 4999.5
$ws.Cells.Item(122, 11).Value = 8987.0772
$ws.Cells.Item(122, 12).Value = 14998.5
$ws.Cells.Item(122, 13).Value = -6537.0772
$ws.Cells.Item(122, 14).Value = -19898.5

# hunk 22: sheet LTW
$ws.Cells.Item(126, 8).Value = 41372.332
$ws.Cells.Item(126, 9).Value = 63853.145
$ws.Cells.Item(126, 11).Value = 191559.435
$ws.Cells.Item(126, 13).Value = -189089.435

# hunk 23: sheet LTW
$ws.Cells.Item(136, 8).Value = 6919.364
$ws.Cells.Item(136, 9).Value = 2329
$ws.Cells.Item(136, 10).Value = 12427.8
$ws.Cells.Item(136, 11).Value = 6987
$ws.Cells.Item(136, 12).Value = 37283.39999999999
$ws.Cells.Item(136, 13).Value = -4437
$ws.Cells.Item(136, 14).Value = -42383.39999999999

$ws = $wb.Worksheets.Item("WVR")
# hunk 24: sheet WVR
$ws.Cells.Item(39, 8).Value = 25000
$ws.Cells.Item(39, 10).Value = 25000
$ws.Cells.Item(39, 12).Value = 25000
$ws.Cells.Item(39, 14).Value = -25826

# hunk 25: sheet WVR
$ws.Cells.Item(42, 8).Value = 25022
$ws.Cells.Item(42, 10).Value = 20000
$ws.Cells.Item(42, 12).Value = 20000
$ws.Cells.Item(42, 14).Value = -20756

# hunk 26: sheet WVR
$ws.Cells.Item(100, 8).Value = 31879.85
$ws.Cells.Item(100, 9).Value = 20162.438
$ws.Cells.Item(100, 10).Value = 78749.5
$ws.Cells.Item(100, 11).Value = 40324.876
$ws.Cells.Item(100, 12).Value = 157499
$ws.Cells.Item(100, 13).Value = -39783.876
$ws.Cells.Item(100, 14).Value = -158581

# hunk 27: sheet WVR
$ws.Cells.Item(122, 8).Value = 17666.055
$ws.Cells.Item(122, 9).Value = 3248.0356
$ws.Cells.Item(122, 11).Value = 9744.106800000001
$ws.Cells.Item(122, 13).Value = -7294.106800000001

# hunk 28: sheet WVR
$ws.Cells.Item(126, 8).Value = 24563.947
$ws.Cells.Item(126, 9).Value = 37407.453
$ws.Cells.Item(126, 11).Value = 112222.359
$ws.Cells.Item(126, 13).Value = -109752.359

# hunk 29: sheet WVR
$ws.Cells.Item(132, 8).Value = 9236.415000000001
$ws.Cells.Item(132, 9).Value = 10334.5
$ws.Cells.Item(132, 10).Value = 5857.6924
$ws.Cells.Item(132, 11).Value = 31003.5
$ws.Cells.Item(132, 12).Value = 17573.0772
$ws.Cells.Item(132, 13).Value = -28473.5
$ws.Cells.Item(132, 14).Value = -22633.0772

# hunk 30: sheet WVR
$ws.Cells.Item(135, 8).Value = 8443376
$ws.Cells.Item(135, 10).Value = 8443376
$ws.Cells.Item(135, 12).Value = 8443376
$ws.Cells.Item(135, 14).Value = -8453516

# hunk 31: sheet WVR
$ws.Cells.Item(136, 8).Value = 2769.9167
$ws.Cells.Item(136, 9).Value = 1178.4286
$ws.Cells.Item(136, 10).Value = 4998
$ws.Cells.Item(136, 11).Value = 3535.2858
$ws.Cells.Item(136, 12).Value = 14994
$ws.Cells.Item(136, 13).Value = -985.2857999999997
$ws.Cells.Item(136, 14).Value = -20094
